$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("monoplane_spar_layup")
$ws.Range("V1").Value = "I_x"
$ws.Range("T1").Value = "I_x (outer)"
$ws.Range("U1").Value = "I_x (inner)"
